# Guild.xlsx update: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet formerly called "Property1" becomes the generic "DataNode"
# sheet, and the sheet formerly called "Record" becomes the generic
# "DataTable" sheet. The author was last working on the DataTable sheet,
# so that is the sheet that ends up active/selected when the file is saved.

$wb = $excel.ActiveWorkbook

$nodeSheet  = $wb.Worksheets.Item(1)   # was "Property1"
$tableSheet = $wb.Worksheets.Item(2)   # was "Record"

$nodeSheet.Name  = "DataNode"
$tableSheet.Name = "DataTable"

# Make the DataTable sheet the active/selected tab (matches the saved
# workbook view: tabSelected moves from sheet 1 to sheet 2).
$tableSheet.Activate()
